$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "Save" in H1, matching the formatting used by the other header cells
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Fill H2:H10 with 0
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
